$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.385.60'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '1.779.67'
$ws.Range("E3").Value = '  +3.74%  '
$c = $ws.Range("D4")
$c.Value = '''1.001'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$c = $ws.Range("D5")
$c.Value = '''313.96'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("E6").Value = '  -0.05%  '
$c = $ws.Range("D7")
$c.Value = '''0.5243'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +9.31%  '
$c = $ws.Range("D8")
$c.Value = '''0.3766'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +9.07%  '
$c = $ws.Range("D10")
$c.Value = '''0.07411'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.25%  '
$ws.Range("E11").Value = '  +5.35%  '
$c = $ws.Range("D12")
$c.Value = '''1.001'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '
$c = $ws.Range("D13")
$c.Value = '''20.70'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.93%  '
$c = $ws.Range("D14")
$c.Value = '''6.106'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +4.60%  '
$ws.Range("D15").Value = '1.779.10'
$ws.Range("E15").Value = '  +3.65%  '
$c = $ws.Range("D16")
$c.Value = '''6.988'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +2.47%  '
$c = $ws.Range("D17")
$c.Value = '''89.66'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.80%  '
$c = $ws.Range("D18")
$c.Value = '''0.00001057'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.41%  '
$c = $ws.Range("D19")
$c.Value = '''0.06437'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.96%  '
$c = $ws.Range("D20")
$c.Value = '''0.9996'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '
$c = $ws.Range("D21")
$c.Value = '''16.81'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.22%  '
$c = $ws.Range("D22")
$c.Value = '''5.890'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +4.73%  '
$ws.Range("D23").Value = '27.429.90'
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("E24").Value = '  +4.40%  '
$c = $ws.Range("D25")
$c.Value = '''2.092'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = '''155.46'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +3.14%  '
$ws.Range("E27").Value = '  +1.44%  '
$c = $ws.Range("D28")
$c.Value = '''2.361'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +14.75%  '
$ws.Range("D29").Value = '1.986.59'
$ws.Range("E29").Value = '  +3.92%  '
$c = $ws.Range("D30")
$c.Value = '''121.16'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.47%  '
$c = $ws.Range("D31")
$c.Value = '''1.090'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +5.10%  '
$c = $ws.Range("D32")
$c.Value = '''0.1016'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +9.97%  '
$c = $ws.Range("D33")
$c.Value = '''5.610'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +5.70%  '
$c = $ws.Range("D34")
$c.Value = '''3.637'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.05%  '
$c = $ws.Range("D35")
$c.Value = '''0.02257'
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = '''0.06000'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.74%  '
$c = $ws.Range("D37")
$c.Value = '''11.32'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.68%  '
$ws.Range("E38").Value = '  +3.46%  '
$c = $ws.Range("D39")
$c.Value = '''4.899'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +4.06%  '
$c = $ws.Range("D40")
$c.Value = '''0.6135'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D41")
$c.Value = '''1.433'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.85%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D42")
$c.Value = '''8.208'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +9.52%  '
$c = $ws.Range("D43")
$c.Value = '''1.133'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +4.69%  '
$c = $ws.Range("D44")
$c.Value = '''13.21'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = '''0.5799'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +4.30%  '
$c = $ws.Range("D46")
$c.Value = '''3.626'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.18%  '
$c = $ws.Range("D47")
$c.Value = '''121.79'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.72%  '
$c = $ws.Range("D48")
$c.Value = '''1.896'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +4.05%  '
$c = $ws.Range("D49")
$c.Value = '''1.118'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.81%  '
$c = $ws.Range("D50")
$c.Value = '''0.06735'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.51%  '
$c = $ws.Range("D51")
$c.Value = '''70.93'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.52%  '
